# Updates to AVL and SOCDtINtY to align sales data.
#
# - SoCDTtiNTY-psgr: LDVs (row 2) share-that-is-new drops from 8.6% to 6.5%
#   across all forecast years (B2:H2).
# - SoCDTtiNTY-frgt: HDVs (row 3) share-that-is-new rises from 3.55% to
#   5.75% across all forecast years (B3:H3).
# - Leave the active selection on each sheet where the author last left it.

$wb = $excel.ActiveWorkbook

# --- SoCDTtiNTY-psgr: LDVs share-that-is-new 8.6% -> 6.5% ---
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsPsgr.Range("B2:H2").Value = 0.065
$wsPsgr.Range("B3:H3").Select() | Out-Null

# --- SoCDTtiNTY-frgt: HDVs share-that-is-new 3.55% -> 5.75% ---
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")
$wsFrgt.Range("B3:H3").Value = 0.0575
$wsFrgt.Range("B3:H3").Select() | Out-Null

# --- About sheet: leave selection where the author last left it ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("L12").Select() | Out-Null
